$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H31").Value = 4500
$ws.Range("I31").Value = 1000
$ws.Range("K31").Value = 3000
$ws.Range("M31").Value = -2770

$ws.Range("H33").Value = 215.3
$ws.Range("I33").Value = 228.11111
$ws.Range("K33").Value = 228.11111
$ws.Range("M33").Value = 0.8888900000000035

$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()

$ws.Range("H62").Value = 2799.4348
$ws.Range("J62").Value = 3606.5454
$ws.Range("L62").Value = 3606.5454
$ws.Range("N62").Value = -4854.5454

$ws.Range("H65").Value = 2799.4348
$ws.Range("J65").Value = 3606.5454
$ws.Range("L65").Value = 18032.727
$ws.Range("N65").Value = -24272.727

$ws.Range("H96").Value = 41667060
$ws.Range("I96").Value = 41667060
$ws.Range("K96").Value = 125001180
$ws.Range("M96").Value = -124999807

$ws.Range("H116").Value = 5660.6
$ws.Range("I116").Value = 2866.6667
$ws.Range("J116").Value = 6858
$ws.Range("K116").Value = 2866.6667
$ws.Range("L116").Value = 6858
$ws.Range("M116").Value = 575.3332999999998
$ws.Range("N116").Value = -13742

$ws.Range("H132").Value = 3144.5356
$ws.Range("I132").Value = 3719.1738
$ws.Range("J132").Value = 501.2
$ws.Range("K132").Value = 11157.5214
$ws.Range("L132").Value = 1503.6
$ws.Range("M132").Value = -8627.5214
$ws.Range("N132").Value = -6563.6

$ws.Range("H135").Value = 27787028
$ws.Range("I135").Value = 883.7692
$ws.Range("K135").Value = 7953.922799999999
$ws.Range("M135").Value = -5418.922799999999

$ws.Range("H141").Value = 1978.9474
$ws.Range("I141").Value = 1371.1538
$ws.Range("K141").Value = 4113.4614
$ws.Range("M141").Value = 1066.5386

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()

$ws.Range("H32").Value = 6283.3877
$ws.Range("I32").Value = 4540.1274
$ws.Range("K32").Value = 4540.1274
$ws.Range("M32").Value = -4253.1274

$ws.Range("H35").Value = 5500
$ws.Range("I35").Value = 5000
$ws.Range("J35").Value = 6000
$ws.Range("K35").Value = 5000
$ws.Range("L35").Value = 6000
$ws.Range("M35").Value = -4594
$ws.Range("N35").Value = -6812

$ws.Range("H45").Value = 2170.6897
$ws.Range("I45").Value = 2810.25
$ws.Range("J45").Value = 1719.2354
$ws.Range("K45").Value = 2810.25
$ws.Range("L45").Value = 1719.2354
$ws.Range("M45").Value = -2433.25
$ws.Range("N45").Value = -2473.2354

$ws.Range("H61").Value = 1666.5946
$ws.Range("I61").Value = 1477.3549
$ws.Range("J61").Value = 2644.3333
$ws.Range("K61").Value = 1477.3549
$ws.Range("L61").Value = 2644.3333
$ws.Range("M61").Value = -1265.3549
$ws.Range("N61").Value = -3068.3333

$ws.Range("H63").Value = 5000
$ws.Range("J63").Value = 5000
$ws.Range("L63").Value = 5000
$ws.Range("N63").Value = -6372

$ws.Range("H66").Value = 5000
$ws.Range("J66").Value = 5000
$ws.Range("L66").Value = 25000
$ws.Range("N66").Value = -31864

$ws.Range("H97").Value = 1152.7142
$ws.Range("I97").Value = 1311.375
$ws.Range("J97").Value = 645
$ws.Range("K97").Value = 1311.375
$ws.Range("L97").Value = 645
$ws.Range("M97").Value = -815.375
$ws.Range("N97").Value = -1637

$ws.Range("H102").Value = 1699.375
$ws.Range("I102").Value = 1423.75
$ws.Range("K102").Value = 1423.75
$ws.Range("M102").Value = 198.25

$ws.Range("H110").Value = 901.3333
$ws.Range("I110").Value = 764
$ws.Range("K110").Value = 764
$ws.Range("M110").Value = 1281

$ws.Range("H122").Value = 3336.125
$ws.Range("I122").Value = 3336.125
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10008.375
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7558.375
$ws.Range("N122").ClearContents()

$ws.Range("H136").Value = 1666.5946
$ws.Range("I136").Value = 1477.3549
$ws.Range("J136").Value = 2644.3333
$ws.Range("K136").Value = 4432.0647
$ws.Range("L136").Value = 7932.999899999999
$ws.Range("M136").Value = -1882.0647
$ws.Range("N136").Value = -13032.9999

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H20").Value = 998.7143
$ws.Range("I20").Value = 1036.75
$ws.Range("J20").Value = 948
$ws.Range("K20").Value = 1036.75
$ws.Range("L20").Value = 948
$ws.Range("M20").Value = -789.75
$ws.Range("N20").Value = -1442

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H50").Value = 13493.333
$ws.Range("J50").Value = 13493.333
$ws.Range("L50").Value = 13493.333
$ws.Range("N50").Value = -14743.333

$ws.Range("H58").Value = 30970.412
$ws.Range("I58").Value = 1727.2
$ws.Range("J58").Value = 72746.42999999999
$ws.Range("K58").Value = 1727.2
$ws.Range("L58").Value = 72746.42999999999
$ws.Range("M58").Value = -1524.2
$ws.Range("N58").Value = -73152.42999999999

$ws.Range("H132").Value = 17625.854
$ws.Range("I132").Value = 27377.1
$ws.Range("J132").Value = 3695.5
$ws.Range("K132").Value = 82131.29999999999
$ws.Range("L132").Value = 11086.5
$ws.Range("M132").Value = -79601.29999999999
$ws.Range("N132").Value = -16146.5

$ws.Range("H136").Value = 30970.412
$ws.Range("I136").Value = 1727.2
$ws.Range("J136").Value = 72746.42999999999
$ws.Range("K136").Value = 5181.6
$ws.Range("L136").Value = 218239.29
$ws.Range("M136").Value = -2631.6
$ws.Range("N136").Value = -223339.29

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H37").Value = 125010000
$ws.Range("J37").Value = 125010000
$ws.Range("L37").Value = 375030000
$ws.Range("N37").Value = -375030224

$ws.Range("H131").Value = 765
$ws.Range("J131").Value = 765
$ws.Range("L131").Value = 2295
$ws.Range("N131").Value = -12375

$ws.Range("H132").Value = 1244.8572
$ws.Range("I132").Value = 483.33334
$ws.Range("J132").Value = 1452.5454
$ws.Range("K132").Value = 4350.00006
$ws.Range("L132").Value = 13072.9086
$ws.Range("M132").Value = -1820.00006
$ws.Range("N132").Value = -18132.9086

$ws.Range("H140").Value = 2132.1667
$ws.Range("I140").Value = 858
$ws.Range("K140").Value = 2574
$ws.Range("M140").Value = 2606

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H58").Value = 16672167
$ws.Range("I58").Value = 3000
$ws.Range("J58").Value = 20006000
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 20006000
$ws.Range("M58").Value = -2723
$ws.Range("N58").Value = -20006554

$ws.Range("H80").Value = 3514.4333
$ws.Range("I80").Value = 3201.8333
$ws.Range("J80").Value = 3722.8333
$ws.Range("K80").Value = 3201.8333
$ws.Range("L80").Value = 3722.8333
$ws.Range("M80").Value = -2203.8333
$ws.Range("N80").Value = -5718.8333

$ws.Range("H83").Value = 3514.4333
$ws.Range("I83").Value = 3201.8333
$ws.Range("J83").Value = 3722.8333
$ws.Range("K83").Value = 16009.1665
$ws.Range("L83").Value = 18614.1665
$ws.Range("M83").Value = -11017.1665
$ws.Range("N83").Value = -28598.1665

$ws.Range("H122").Value = 2002
$ws.Range("I122").Value = 1007
$ws.Range("J122").Value = 2499.5
$ws.Range("K122").Value = 3021
$ws.Range("L122").Value = 7498.5
$ws.Range("M122").Value = -571
$ws.Range("N122").Value = -12398.5

$ws.Range("H132").Value = 31520.422
$ws.Range("I132").Value = 6030
$ws.Range("J132").Value = 86749.664
$ws.Range("K132").Value = 18090
$ws.Range("L132").Value = 260248.992
$ws.Range("M132").Value = -15560
$ws.Range("N132").Value = -265308.992

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H61").Value = 4447.684
$ws.Range("J61").Value = 7287.2856
$ws.Range("L61").Value = 7287.2856
$ws.Range("N61").Value = -7691.2856

$ws.Range("H113").Value = 4447.684
$ws.Range("J113").Value = 7287.2856
$ws.Range("L113").Value = 7287.2856
$ws.Range("N113").Value = -11627.2856

$ws.Range("H122").Value = 1228478.6
$ws.Range("I122").Value = 2181349.8
$ws.Range("J122").Value = 3358.5715
$ws.Range("K122").Value = 6544049.399999999
$ws.Range("L122").Value = 10075.7145
$ws.Range("M122").Value = -6541599.399999999
$ws.Range("N122").Value = -14975.7145

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H100").Value = 299.5
$ws.Range("I100").Value = 299.33334
$ws.Range("J100").Value = 300
$ws.Range("K100").Value = 598.66668
$ws.Range("L100").Value = 600
$ws.Range("M100").Value = -57.66668000000004
$ws.Range("N100").Value = -1682

$ws.Range("H122").Value = 2035.4445
$ws.Range("I122").Value = 1950
$ws.Range("J122").Value = 2462.6667
$ws.Range("K122").Value = 5850
$ws.Range("L122").Value = 7388.000100000001
$ws.Range("M122").Value = -3400
$ws.Range("N122").Value = -12288.0001

$ws.Range("H136").Value = 34484500
$ws.Range("I136").Value = 41668100
$ws.Range("J136").Value = 3220.8
$ws.Range("K136").Value = 125004300
$ws.Range("L136").Value = 9662.400000000001
$ws.Range("M136").Value = -125001750
$ws.Range("N136").Value = -14762.4

Write-Output "Edit applied successfully"
